$wb = $excel.ActiveWorkbook

# Scheduled market-data refresh: update currentAveragePrice / LevePrice / LeveProfit
# columns (H-N) across the per-job Leve sheets. Values sourced from the latest
# Universalis market snapshot; row/column layout is untouched.

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 4758.9
$ws.Range("I19").Value = 2051.0833
$ws.Range("K19").Value = 2051.0833
$ws.Range("M19").Value = -1876.0833
$ws.Range("H40").Value = 3921.1428
$ws.Range("I40").Value = 1689.9
$ws.Range("K40").Value = 1689.9
$ws.Range("M40").Value = -1514.9
$ws.Range("H43").Value = 3622.5
$ws.Range("J43").Value = 1499
$ws.Range("L43").Value = 1499
$ws.Range("N43").Value = -1637
$ws.Range("H70").Value = 2270.8572
$ws.Range("I70").Value = 1998.5
$ws.Range("K70").Value = 5995.5
$ws.Range("M70").Value = -5725.5
$ws.Range("H73").Value = 2270.8572
$ws.Range("I73").Value = 1998.5
$ws.Range("K73").Value = 5995.5
$ws.Range("M73").Value = -5059.5
$ws.Range("H92").Value = 28254.637
$ws.Range("I92").Value = 797.9259
$ws.Range("J92").Value = 151809.83
$ws.Range("K92").Value = 797.9259
$ws.Range("L92").Value = 151809.83
$ws.Range("M92").Value = 450.0741
$ws.Range("N92").Value = -154305.83
$ws.Range("H98").Value = 620.5833
$ws.Range("I98").Value = 552
$ws.Range("K98").Value = 552
$ws.Range("M98").Value = 946
$ws.Range("H113").Value = 40570.145
$ws.Range("I113").Value = 59499.445
$ws.Range("K113").Value = 59499.445
$ws.Range("M113").Value = -56245.445
$ws.Range("H122").Value = 620.5833
$ws.Range("I122").Value = 552
$ws.Range("K122").Value = 1656
$ws.Range("M122").Value = 794
$ws.Range("H132").Value = 1114079.8
$ws.Range("I132").Value = 1114079.8
$ws.Range("K132").Value = 3342239.4
$ws.Range("M132").Value = -3339709.4
$ws.Range("H138").Value = 1401.5
$ws.Range("I138").Value = 911.2692
$ws.Range("K138").Value = 2733.8076
$ws.Range("M138").Value = 2406.1924

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1760.9
$ws.Range("J2").Value = 1474.8
$ws.Range("L2").Value = 1474.8
$ws.Range("N2").Value = -1700.8
$ws.Range("H45").Value = 2122.24
$ws.Range("I45").Value = 1217.5333
$ws.Range("K45").Value = 1217.5333
$ws.Range("M45").Value = -840.5333000000001
$ws.Range("H102").Value = 1804.0435
$ws.Range("I102").Value = 1808.7727
$ws.Range("K102").Value = 1808.7727
$ws.Range("M102").Value = -186.7727
$ws.Range("H116").Value = 1760.9
$ws.Range("J116").Value = 1474.8
$ws.Range("L116").Value = 1474.8
$ws.Range("N116").Value = -6062.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1760.9
$ws.Range("J3").Value = 1474.8
$ws.Range("L3").Value = 1474.8
$ws.Range("N3").Value = -1702.8
$ws.Range("H86").Value = 41667624
$ws.Range("J86").Value = 1000.0909
$ws.Range("L86").Value = 1000.0909
$ws.Range("N86").Value = -3246.0909
$ws.Range("H89").Value = 41667624
$ws.Range("J89").Value = 1000.0909
$ws.Range("L89").Value = 5000.4545
$ws.Range("N89").Value = -16232.4545
$ws.Range("H94").Value = 5000924
$ws.Range("I94").Value = 6250738.5
$ws.Range("K94").Value = 6250738.5
$ws.Range("M94").Value = -6250287.5
$ws.Range("H105").Value = 55557068
$ws.Range("I105").Value = 83334690
$ws.Range("K105").Value = 83334690
$ws.Range("M105").Value = -83332943

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 76665.86
$ws.Range("I62").Value = 105732.6
$ws.Range("K62").Value = 105732.6
$ws.Range("M62").Value = -105108.6
$ws.Range("H65").Value = 76665.86
$ws.Range("I65").Value = 105732.6
$ws.Range("K65").Value = 528663
$ws.Range("M65").Value = -525543
$ws.Range("H107").Value = 1242.9429
$ws.Range("I107").Value = 1015.85187
$ws.Range("K107").Value = 1015.85187
$ws.Range("M107").Value = 904.14813
$ws.Range("H122").Value = 5532.875
$ws.Range("I122").Value = 5510.8125
$ws.Range("J122").Value = 5577
$ws.Range("K122").Value = 16532.4375
$ws.Range("L122").Value = 16731
$ws.Range("M122").Value = -14082.4375
$ws.Range("N122").Value = -21631

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H76").Value = 6814
$ws.Range("J76").Value = 6814
$ws.Range("L76").Value = 20442
$ws.Range("N76").Value = -21208
$ws.Range("H79").Value = 6814
$ws.Range("J79").Value = 6814
$ws.Range("L79").Value = 20442
$ws.Range("N79").Value = -23094
$ws.Range("H92").Value = 186.90909
$ws.Range("I92").Value = 95.25
$ws.Range("J92").Value = 239.28572
$ws.Range("K92").Value = 285.75
$ws.Range("L92").Value = 717.85716
$ws.Range("M92").Value = 962.25
$ws.Range("N92").Value = -3213.85716
$ws.Range("H95").Value = 15027
$ws.Range("J95").Value = 15027
$ws.Range("L95").Value = 45081
$ws.Range("N95").Value = -49199
$ws.Range("H97").Value = 152
$ws.Range("H102").Value = 2999.5
$ws.Range("I102").Value = 2999.5
$ws.Range("K102").Value = 8998.5
$ws.Range("M102").Value = -6564.5
$ws.Range("H111").Value = 14744
$ws.Range("I111").Value = 0
$ws.Range("K111").Value = 0
$ws.Range("H120").Value = 2399.5
$ws.Range("I120").Value = 2399.5
$ws.Range("K120").Value = 7198.5
$ws.Range("M120").Value = -2360.5
$ws.Range("H132").Value = 1771.5454
$ws.Range("I132").Value = 1622
$ws.Range("J132").Value = 1857
$ws.Range("K132").Value = 14598
$ws.Range("L132").Value = 16713
$ws.Range("M132").Value = -12068
$ws.Range("N132").Value = -21773
$ws.Range("H138").Value = 4061.3333
$ws.Range("I138").Value = 2885.0908
$ws.Range("K138").Value = 8655.2724
$ws.Range("M138").Value = -3515.2724

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3298.7144
$ws.Range("I80").Value = 2400
$ws.Range("J80").Value = 3658.2
$ws.Range("K80").Value = 2400
$ws.Range("L80").Value = 3658.2
$ws.Range("M80").Value = -1402
$ws.Range("N80").Value = -5654.2
$ws.Range("H83").Value = 3298.7144
$ws.Range("I83").Value = 2400
$ws.Range("J83").Value = 3658.2
$ws.Range("K83").Value = 12000
$ws.Range("L83").Value = 18291
$ws.Range("M83").Value = -7008
$ws.Range("N83").Value = -28275
$ws.Range("H97").Value = 1245.6
$ws.Range("I97").Value = 1109.2759
$ws.Range("K97").Value = 1109.2759
$ws.Range("M97").Value = -613.2759000000001
$ws.Range("H117").Value = 58000
$ws.Range("J117").Value = 58000
$ws.Range("L117").Value = 58000
$ws.Range("N117").Value = -64884
$ws.Range("H126").Value = 52635440
$ws.Range("J126").Value = 5441.143
$ws.Range("L126").Value = 16323.429
$ws.Range("N126").Value = -21263.429

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4196.6665
$ws.Range("J40").Value = 4344.25
$ws.Range("L40").Value = 4344.25
$ws.Range("N40").Value = -4616.25
$ws.Range("H46").Value = 2820.5557
$ws.Range("J46").Value = 3284.8572
$ws.Range("L46").Value = 3284.8572
$ws.Range("N46").Value = -3660.8572
$ws.Range("H82").Value = 1696.0555
$ws.Range("I82").Value = 2034.2
$ws.Range("K82").Value = 2034.2
$ws.Range("M82").Value = -1673.2
$ws.Range("H85").Value = 1696.0555
$ws.Range("I85").Value = 2034.2
$ws.Range("K85").Value = 2034.2
$ws.Range("M85").Value = -786.2
$ws.Range("H100").Value = 4748.3076
$ws.Range("J100").Value = 4407
$ws.Range("L100").Value = 4407
$ws.Range("N100").Value = -5489
$ws.Range("H122").Value = 4395.316
$ws.Range("I122").Value = 3347.3076
$ws.Range("K122").Value = 10041.9228
$ws.Range("M122").Value = -7591.9228

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3516.8333
$ws.Range("I122").Value = 3146.1052
$ws.Range("J122").Value = 4925.6333
$ws.Range("K122").Value = 9438.3156
$ws.Range("L122").Value = 14776.8
$ws.Range("M122").Value = -6988.3156
$ws.Range("N122").Value = -19676.8

# Row 111 on CUL: the refreshed snapshot has no LeveProfitNQ figure for this
# leve, so the M111 cell is cleared entirely rather than re-valued.
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("M111").ClearContents()

Write-Output "Applied 201 cell updates + 1 cell clear across 8 sheets"